# "san sefiro fungerar nu"
# Swap the document's Roboto font usages over to "San Sefiro" everywhere a
# w:rFonts element already specifies Roboto (ascii/hAnsi/eastAsia). Runs and
# paragraph marks that currently have no font set are left untouched so we
# don't introduce new rFonts/rPr blocks that weren't there before.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    if ($rng.Font.NameAscii -eq "Roboto") {
        $rng.Font.NameAscii   = "San Sefiro"
        $rng.Font.NameOther   = "San Sefiro"
        $rng.Font.NameFarEast = "San Sefiro"
    }
}
